$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. row195 / B195 time value changed (0.979... -> 0.479...) ---
$ws.Range("B195").Value = 0.47916666666666669

# --- 2. row238: chargebacks turned in to GI -> C238/D238 updated ---
$ws.Range("C238").Value = "PTR Chores"
$ws.Range("D238").Value = "chargebacks"

# --- 3. drop the very last (fully blank) row so the sheet ends at row 254 ---
$ws.Rows.Item(255).Delete()

# --- 4. fill in the 7/21-7/22 entries that used to be blank placeholder rows ---
# Make sure each new data row carries the same date/time number formats (styles)
# as the existing data rows by copying from row 238 (s="1" date, s="2" time) first.
$ws.Range("A238:B238").Copy($ws.Range("A239:B239"))
$ws.Range("A238:B238").Copy($ws.Range("A240:B240"))
$ws.Range("A238:B238").Copy($ws.Range("A241:B241"))
$ws.Range("A238:B238").Copy($ws.Range("A242:B242"))
$ws.Range("A238:B238").Copy($ws.Range("A243:B243"))
$ws.Range("A238:B238").Copy($ws.Range("A244:B244"))
$ws.Range("A238:B238").Copy($ws.Range("A245:B245"))
$ws.Range("A238:B238").Copy($ws.Range("A246:B246"))
$ws.Range("A238:B238").Copy($ws.Range("A247:B247"))
$ws.Range("A238:B238").Copy($ws.Range("A248:B248"))
$ws.Range("A238:B238").Copy($ws.Range("A249:B249"))
$ws.Range("A238:B238").Copy($ws.Range("A250:B250"))

# Row 239 - 7/21/16
$ws.Range("A239").Value = 42572
$ws.Range("B239").Value = 0.67638888888888893
$ws.Range("C239").Value = "PTR Chores"
$ws.Range("D239").Value = "PTR-IHC Smartsheet updates"

# Row 240 - 7/21/16
$ws.Range("A240").Value = 42572
$ws.Range("B240").Value = 0.68402777777777779
$ws.Range("C240").Value = "Other Imaging"
$ws.Range("D240").Value = "confocal tables updates"

# Row 241 - 7/21/16
$ws.Range("A241").Value = 42572
$ws.Range("B241").Value = 0.69444444444444453
$ws.Range("C241").Value = "KK - Post Processing"
$ws.Range("D241").Value = "slidemaps"
$ws.Range("E241").Value = "20160527_100_KK"

# Row 242 - 7/21/16
$ws.Range("A242").Value = 42572
$ws.Range("B242").Value = 0.70486111111111116
$ws.Range("C242").Value = "KK - Post Processing"
$ws.Range("D242").Value = "slidemaps"
$ws.Range("E242").Value = "20160602_100_KK"

# Row 243 - 7/21/16
$ws.Range("A243").Value = 42572
$ws.Range("B243").Value = 0.71527777777777779
$ws.Range("C243").Value = "SH - Post Process"
$ws.Range("D243").Value = "slidemaps"
$ws.Range("E243").Value = "20160712-03_SH"

# Row 244 - 7/21/16
$ws.Range("A244").Value = 42572
$ws.Range("B244").Value = 0.72916666666666663
$ws.Range("C244").Value = "depart"

# Row 245 - 7/22/16 (guesses)
$ws.Range("A245").Value = 42573
$ws.Range("B245").Value = 0.35416666666666669
$ws.Range("C245").Value = "PTO/ STO"
$ws.Range("D245").Value = "STO"

# Row 246 - 7/22/16 (guesses)
$ws.Range("A246").Value = 42573
$ws.Range("B246").Value = 0.48958333333333331
$ws.Range("C246").Value = "PTR Chores"
$ws.Range("D246").Value = "chargebacks"

# Row 247 - 7/22/16 (guesses)
$ws.Range("A247").Value = 42573
$ws.Range("B247").Value = 0.54166666666666663
$ws.Range("C247").Value = "Reiser Lab The Box"
$ws.Range("D247").Value = "analysis"

# Row 248 - 7/22/16 (guesses)
$ws.Range("A248").Value = 42573
$ws.Range("B248").Value = 0.5625
$ws.Range("C248").Value = "XZ - Post Processing"

# Row 249 - 7/22/16 (guesses)
$ws.Range("A249").Value = 42573
$ws.Range("B249").Value = 0.60416666666666663
$ws.Range("C249").Value = "SH - Post Process"

# Row 250 - 7/22/16 (guesses)
$ws.Range("A250").Value = 42573
$ws.Range("B250").Value = 0.70833333333333337
$ws.Range("C250").Value = "depart"

# --- 5. move the selection/scroll position to reflect where the log now ends ---
$ws.Range("A251:XFD251").Select()
